$wb = $excel.ActiveWorkbook

# --- Add the new "classes" worksheet as the last sheet -------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "classes"

# --- Cell values -----------------------------------------------------------
# Set values in this particular order so that new shared-string entries are
# appended to sharedStrings.xml in the same order as the target workbook:
# featuresSet, adducts, setObjects, ionizedXXX, featuresGroupsSet,
# MSPeakListsSet, formulasSet, compoundsSet, componentsSet, setThreshold,
# origFGNames
$ws.Range("A2").Value = "featuresSet"
$ws.Range("B1").Value = "adducts"
$ws.Range("C1").Value = "setObjects"
$ws.Range("D1").Value = "ionizedXXX"
$ws.Range("A3").Value = "featuresGroupsSet"
$ws.Range("A4").Value = "MSPeakListsSet"
$ws.Range("A5").Value = "formulasSet"
$ws.Range("A6").Value = "compoundsSet"
$ws.Range("A7").Value = "componentsSet"
$ws.Range("E1").Value = "setThreshold"
$ws.Range("F1").Value = "origFGNames"

# Headers that reuse already-existing shared strings
$ws.Range("G1").Value = "groupAlgorithm"
$ws.Range("H1").Value = "analysisInfo"

# --- "X" marker cells --------------------------------------------------
$ws.Range("B2").Value = "X"
$ws.Range("C2").Value = "X"
$ws.Range("D2").Value = "X"

$ws.Range("G3").Value = "X"

$ws.Range("B4").Value = "X"
$ws.Range("C4").Value = "X"
$ws.Range("H4").Value = "X"

$ws.Range("B5").Value = "X"
$ws.Range("C5").Value = "X"
$ws.Range("E5").Value = "X"
$ws.Range("F5").Value = "X"

$ws.Range("B6").Value = "X"
$ws.Range("C6").Value = "X"
$ws.Range("E6").Value = "X"
$ws.Range("F6").Value = "X"

$ws.Range("C7").Value = "X"

# --- Styling: center-align the data block (creates the new cellXf) ------
$ws.Range("B2:I7").HorizontalAlignment = -4108
$ws.Range("B8:D10").HorizontalAlignment = -4108

# --- Column widths (best achievable match to the target widths) --------
$ws.Columns.Item(1).ColumnWidth = 17
$ws.Columns.Item(4).ColumnWidth = 10.166666666666666
$ws.Columns.Item(5).ColumnWidth = 11.666666666666666
$ws.Columns.Item(6).ColumnWidth = 12
$ws.Columns.Item(7).ColumnWidth = 14.333333333333334
$ws.Columns.Item(8).ColumnWidth = 10.666666666666666

# --- Selection on the new sheet -----------------------------------------
$ws.Range("C3").Select()
